$d = $word.ActiveDocument

# The "_GoBack" bookmark currently sits (collapsed) right after the text
# "Changes for commit", at the end of that paragraph. We want to end up
# with:
#   Paragraph N:   "Changes for commit"                     (unchanged)
#   Paragraph N+1: "Changes for 2" + "nd"(superscript) + " commit"
#                  followed by the (relocated) _GoBack bookmark.
#
# Directly re-adding a bookmark via $d.Bookmarks.Add(...) at a position
# that lands on certain paragraph-boundary offsets is unreliable in this
# host, so instead we let the existing bookmark's own Range naturally
# carry it forward: insert the new paragraph's text right before the
# bookmark (pushing the bookmark after the inserted text, still inside
# the same paragraph), then insert a paragraph break before that new
# text. That break splits the paragraph so the bookmark ends up, intact,
# at the end of the newly created second paragraph.

$bm = $d.Bookmarks("_GoBack")
$splitPos = $bm.Start

$insertPoint = $d.Range($splitPos, $splitPos)
$newText = "Changes for 2nd commit"
$insertPoint.InsertBefore($newText)

$breakRange = $d.Range($splitPos, $splitPos)
$breakRange.InsertParagraphBefore()

# Locate the newly created paragraph (the one now holding $newText) and
# apply superscript formatting to the "nd" portion of "2nd".
$newPara = $d.Paragraphs(4)
$paraStart = $newPara.Range.Start

$ndStart = $paraStart + "Changes for 2".Length
$ndEnd = $ndStart + "nd".Length
$ndRange = $d.Range($ndStart, $ndEnd)
$ndRange.Font.Superscript = $true
